$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Order date moves forward to reflect the rerouted request
$ws.Range("A2").Value = "06-01-2022"

# New OverageID issued for the rerouted request
$ws.Range("L2").Value = "59071955"
